$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cf = $ws.Range("C53").FormatConditions
Write-Host "Count:" $cf.Count
for ($i=1; $i -le $cf.Count; $i++) {
  $item = $cf.Item($i)
  Write-Host $i $item.Type $item.Operator $item.Formula1
}
